$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 12, 4, 150, 4, 574200, 0.6767144892144892),
    @(1, 12, 4, 200, 4, 574200, 0.6767144892144892),
    @(1, 12, 4, 250, 4, 574200, 0.6767144892144892)
)

$startRow = 9
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
